# CE342_A.xlsx - computation of new "Resenha Novos Clássicos" grades.
# Column F used to hold student e-mails. A new deadline/date column is
# inserted at F, a new grade column "Resenha Novos Clássicos" is inserted
# at G, and the e-mail column shifts to H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Force F1 to stay plain text (it looks like a date, Excel would
# otherwise auto-convert it to a date serial) by setting it as "@"
# (Text) formatted before assigning, then re-apply the same header
# look-and-feel (bold / centered / bordered) that the rest of row 1
# already has, by copying the formatting from an existing header cell
# onto F1, G1 and H1 (keeps them on the shared header style).
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "2020-11-09"
$ws.Range("G1").Value = "Resenha Novos Clássicos"
$ws.Range("H1").Value = "Email"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Data rows (2..44) -------------------------------------------------
# New numeric values for columns F ("Resenha Novos Clássicos" score-ish
# metric) and G (count-like grade) per the commit's computed results.
$fVals = @(11.3, 14.78, 94.78, 0, 0, 13.04, 31.3, 12.17, 18.26, 12.17, 95.65, $null, 0, 0, 0, 93.04, 87.83, 2.61, 89.57, 96.52, 19.13, 15.65, 40.87, 81.74, 54.78, 6.09, 56.52, 7.83, 8.7, 42.61, 0, 27.83, 7.83, 0, 89.57, $null, 92.17, 79.13, 7.83, 93.91, $null, 99.13, 13.91)
$gVals = @(5, 10, 5, $null, 10, 5, 10, 7, 5, 5, 10, 0, $null, $null, $null, 10, 7, 10, 0, 5, 0, 5, 7, 7, 5, 10, 5, 5, 7, 5, 0, 5, 10, 0, 3, 0, 5, 5, 7, 5, 10, 3, 5)

for ($i = 0; $i -lt $fVals.Length; $i++) {
    $row = $i + 2

    # The e-mail currently lives in column F; move it to the new H column
    # before F gets overwritten with the numeric value.
    $email = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 8).Value = $email

    if ($null -ne $fVals[$i]) {
        $ws.Cells.Item($row, 6).Value = $fVals[$i]
    } else {
        $ws.Cells.Item($row, 6).Value = $null
    }

    if ($null -ne $gVals[$i]) {
        $ws.Cells.Item($row, 7).Value = $gVals[$i]
    } else {
        $ws.Cells.Item($row, 7).Value = $null
    }
}
